$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add RSI values for rows 2 and 3 (column E)
$ws.Range("E2").Value = 33.5
$ws.Range("E3").Value = 18.7

# Update 최종점수 (final score) values for rows 2 and 3 (column N)
$ws.Range("N2").Value = 85.87127175646313
$ws.Range("N3").Value = 85.87127175646313
